$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 276; this shifts existing rows 276-311
# down to 277-312, preserving all of their data untouched.
$ws.Rows.Item(276).Insert()

# Populate the newly inserted row 276 with the new weekly price record.
$ws.Range("A276").Value = 9
$ws.Range("B276").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C276").Value = "Metropolitana"
$ws.Range("D276").Value = 44918
$ws.Range("E276").Value = 13
$ws.Range("F276").Value = 100112026
$ws.Range("G276").Value = "Haba"
$ws.Range("H276").Value = "Sin especificar"
$ws.Range("I276").Value = "Primera"
$ws.Range("J276").Value = 70
$ws.Range("K276").Value = 16000
$ws.Range("L276").Value = 17000
$ws.Range("M276").Value = 16500
$ws.Range("N276").Value = "$/saco 25 kilos"
$ws.Range("O276").Value = "Carahue"
$ws.Range("P276").Value = 660
$ws.Range("Q276").Value = 25
$ws.Range("R276").Value = "Hortaliza"
